$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Part 1: the "_GoBack" bookmark currently sits inside a math zone
# (m:e of an m:sSub for "Q_2"). Relocate it out of there by first
# stripping it from that paragraph, then re-adding it at the correct
# spot in the plain-text paragraph below (part 2).
# -----------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $pxml = $p.Range.XML(1)
    if ($pxml -like "*_GoBack*") {
        # Drop the bookmark tags themselves.
        $pxml = $pxml.Replace('<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>', '')

        # Round-tripping a paragraph through Range.XML()/InsertXML() makes
        # this engine re-express <w:numPr> as the legacy <w:listPr><w:ilfo/>
        # form, and stamp a w14:paraId/textId pair that wasn't present
        # before. Undo both so the paragraph is otherwise left untouched.
        $pxml = $pxml.Replace('<w:listPr><w:ilfo w:val="1"/></w:listPr>', '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr>')
        $pxml = [regex]::Replace($pxml, 'w14:paraId="[0-9A-Fa-f]+"\s+w14:textId="[0-9A-Fa-f]+"\s+', '')

        $p.Range.InsertXML($pxml) | Out-Null
        break
    }
}

# -----------------------------------------------------------------------
# Part 2: split the A0-16 paragraph's run right after "...downlink
# frequency of " and insert the (now free) "_GoBack" bookmark between
# the two halves, matching the target diff.
# -----------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("At a downlink frequency of ", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0) | Out-Null
$d.Bookmarks.Add("_GoBack", $rng) | Out-Null
